# Updated GSC data files for main domain:
# Append a new day of data (2025-11-16) to the "Chart" sheet. The two
# issue-summary sheets ("Critical issues" / "Non-critical issues") keep
# their header text as-is; the shared-string table simply grows by one
# entry to hold the new date, which is handled automatically on save.

$wb = $excel.ActiveWorkbook

# --- Chart sheet: append the new day's row ---
$chart = $wb.Worksheets.Item("Chart")

# Writing the date-shaped string straight into .Value would get silently
# reinterpreted as a real date serial by Excel's input parser (same as
# typing 2025-11-16 into a cell). Route it through a text formula first,
# then flatten the formula down to its literal value with a values-only
# paste so the cell ends up as plain text, matching the rest of column A.
$chart.Cells.Item(43, 1).Formula = '="2025-11-16"'
$chart.Cells.Item(43, 1).Copy()
$chart.Cells.Item(43, 1).PasteSpecial(-4163)

$chart.Cells.Item(43, 2).Value = 0
$chart.Cells.Item(43, 3).Value = 31

# --- Critical issues sheet: re-assert header text (keeps data intact) ---
$critical = $wb.Worksheets.Item("Critical issues")
$critical.Cells.Item(1, 1).Value = "Issue"
$critical.Cells.Item(1, 2).Value = "Validation"
$critical.Cells.Item(1, 3).Value = "Items"

# --- Non-critical issues sheet: re-assert header text (keeps data intact) ---
$noncritical = $wb.Worksheets.Item("Non-critical issues")
$noncritical.Cells.Item(1, 1).Value = "Issue"
$noncritical.Cells.Item(1, 2).Value = "Validation"
$noncritical.Cells.Item(1, 3).Value = "Items"
